# This script applies the "overworld enemy data" update described in the
# commit message: the scratch encounter-list notes that had been typed into
# column J (rows 1-32) are removed, and the real Name/GP/Experience data for
# the remaining overworld enemies is filled into columns E:G (rows 23-32) of
# Table2. Table2 is then resized to drop the now-unused trailing row, and the
# active selection is moved to J20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old scratch notes that lived in column J (rows 1-32).
$ws.Range("J1:J32").ClearContents()

# Fill in the real overworld enemy Name/GP/Experience rows that used to be
# missing from Table2 (E23:G32).
$data = @(
  @(23, "StingRat", 220, 1210),
  @(24, "Treant",   150, 1000),
  @(25, "Panther",  255, 830),
  @(26, "Cannibal", 220, 960),
  @(27, "Python",   225, 760),
  @(28, "Roc",      150, 1410),
  @(29, "Roc Baby",  85, 1010),
  @(30, "HugeCell", 255, 1510),
  @(31, "FlameDog", 245, 1720),
  @(32, "BlackLiz",  45, 1300)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("E$r").Value = $row[1]
    $ws.Range("F$r").Value = $row[2]
    $ws.Range("G$r").Value = $row[3]
}

# Table2 no longer needs its trailing blank row now that it ends at row 36.
$lo = $ws.ListObjects.Item("Table2")
$lo.Resize($ws.Range("E1:G36"))

# Match the author's final selection.
$ws.Range("J20").Select()
